$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C120").Value = "combine_2D_plots_v2"
$ws.Range("C121").Value = "add_ASI_background_to_hdf5"
$ws.Range("C122").Value = "create_thm_hdf5"
$ws.Range("C123").Value = "find_irbem_magneticFieldStr"
$ws.Range("C124").Value = "find_irbem_magneticFieldNo"

$ws.Range("C125").Value = "add_thm_hdf5"
$ws.Range("D125").Value = "probeName, outputH5FileStr, omniH5FileStr"

$ws.Range("C126").Value = "write_sc_to_hdf5"
$ws.Range("E126").Value = "Write spacecraft state data to hdf5"
$ws.Range("D126").Value = "h5OutputFile,probeName,time,XYZ_GEO,magFieldStr,NFoot,Lm"

$ws.Range("C127").Value = "geopack_find_magequator"
$ws.Range("C128").Value = "geopack_find_footpoint"
$ws.Range("C129").Value = "geopack_find_curvature"

$ws.Range("C130").Select() | Out-Null
